# Resultados buenos - update order number to numeric value and refresh confidence scores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 was the shared string "4501168528" (text); change it to the real number 4501168528
$ws.Range("C2").Value = 4501168528

# Confidence column (I2:I4) goes from 89.502 to 91.273 (higher confidence)
$ws.Range("I2").Value = 91.273
$ws.Range("I3").Value = 91.273
$ws.Range("I4").Value = 91.273

# Update the active selection / view to match the saved state
[void]$ws.Range("H14").Select()
